# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-01 10:21:50
#
# Applies:
#   - Reordering of the email lists in column G ("Recorded By") for a number
#     of rows (content is the same set of addresses, just re-ordered).
#   - Small numeric corrections in the per-row summary / statistics block
#     (columns L, P, Q).
#   - Row 99 (B3 / ANATOMY / session 2) flips from "Pending" (not yet due)
#     to "Not Recorded" (session date passed with nothing recorded), which
#     also changes its row formatting from the yellow "Pending" style to
#     the pink/red "Not Recorded" style used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column G ("Recorded By") reorderings
# ---------------------------------------------------------------------
$gUpdates = @{
    "G2"   = "Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg"
    "G3"   = "mennatulla.medhat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
    "G7"   = "Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
    "G13"  = "ManaratAleslam-Zefan@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
    "G18"  = "Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg"
    "G19"  = "Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    "G23"  = "Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
    "G29"  = "ManaratAleslam-Zefan@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
    "G34"  = "mennatulla.medhat@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
    "G35"  = "heba@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
    "G50"  = "mennatulla.medhat@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
    "G51"  = "heba@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
    "G66"  = "mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
    "G67"  = "mennatulla.medhat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
    "G77"  = "user@user.com, nourhan.mostafa@med.asu.edu.eg"
    "G82"  = "mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
    "G83"  = "Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    "G93"  = "user@user.com, nourhan.mostafa@med.asu.edu.eg"
    "G98"  = "Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg"
    "G103" = "Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
    "G114" = "Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg"
    "G119" = "Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
}

foreach ($addr in $gUpdates.Keys) {
    $ws.Range($addr).Value = $gUpdates[$addr]
}

# ---------------------------------------------------------------------
# Numeric corrections in the statistics columns
# ---------------------------------------------------------------------
$ws.Range("L7").Value = 6     # Missing Sessions (Group A1)      5 -> 6
$ws.Range("L8").Value = 97    # Pending Sessions (Group A1)     98 -> 97
$ws.Range("P21").Value = 2    # Group B3 row - "Not Recorded" count   1 -> 2
$ws.Range("Q21").Value = 12   # Group B3 row - "Pending" count       13 -> 12

# ---------------------------------------------------------------------
# Row 99: B3 / ANATOMY / session 2 switches from "Pending" to
# "Not Recorded" (its scheduled date has now passed with no attendance
# recorded). Update the status text and re-color the row to match the
# other "Not Recorded" rows (e.g. row 39): pink/red fill, black text.
# ---------------------------------------------------------------------
$ws.Range("I99").Value = "Not Recorded"

$notRecordedRange = $ws.Range("A99:I99")
$notRecordedRange.Interior.Pattern = 1       # xlSolid
$notRecordedRange.Interior.Color = 12695295  # RGB(255,182,193) -> 0x00C1B6FF (BGR)  "#FFB6C1"
$notRecordedRange.Font.Color = 0             # black
$notRecordedRange.HorizontalAlignment = -4108  # xlCenter
$notRecordedRange.VerticalAlignment = -4108    # xlCenter
